# Update the workbook:
#  1. Change every "Förändrad" date (column C, rows 2..490) from 45182 to 45184.
#  2. Append a new row (491) with a new cutting notification record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bump column C (Förändrad) for all data rows from 45182 -> 45184.
for ($r = 2; $r -le 490; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value = 45184
    }
}

# 2. Give row 490 an explicit row height (matches the diff's customHeight flag).
$ws.Rows.Item(490).RowHeight = 15

# 3. Append the new record as row 491.
$newRow = 491
$ws.Cells.Item($newRow, 1).Value = "A 43082-2023"
$ws.Cells.Item($newRow, 2).Value = 45182
$ws.Cells.Item($newRow, 3).Value = 45184
$ws.Cells.Item($newRow, 4).Value = "ÖSTERGÖTLANDS LÄN"
$ws.Cells.Item($newRow, 5).Value = "ÅTVIDABERG"
$ws.Cells.Item($newRow, 7).Value = 2.4
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 0
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = 0
$ws.Cells.Item($newRow, 15).Value = 0
$ws.Cells.Item($newRow, 16).Value = 0
$ws.Cells.Item($newRow, 17).Value = 0

# Match the date-format style already used in columns B and C.
$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"

# R491 keeps the wrap-text style used throughout column R; leave it blank.
$ws.Cells.Item($newRow, 18).WrapText = $true
